# Update "想去人数" (F) and "最低票价" (G) figures on the "展览" and "全部类型"
# sheets to match newly generated site output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 117
$ws1.Range("F3").Value = 5615
$ws1.Range("F6").Value = 930
$ws1.Range("F7").Value = 151
$ws1.Range("F8").Value = 2513
$ws1.Range("F10").Value = 141
$ws1.Range("G10").Value = 65
$ws1.Range("F11").Value = 8
$ws1.Range("F12").Value = 78
$ws1.Range("F14").Value = 2352
$ws1.Range("F15").Value = 330

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 117
$ws4.Range("F3").Value = 5615
$ws4.Range("F8").Value = 930
$ws4.Range("F9").Value = 151
$ws4.Range("F10").Value = 2513
$ws4.Range("F11").Value = 83
$ws4.Range("F12").Value = 141
$ws4.Range("G12").Value = 65
$ws4.Range("F13").Value = 8
$ws4.Range("F15").Value = 79
$ws4.Range("F17").Value = 2352
$ws4.Range("F18").Value = 330
